# Update with restock suggestion
# - Sheet "Forecast Comparison": fill in Week_Start_Date (col B), change
#   "Sales Trend" (col O) from "Stable (-)" to "Stable", update
#   "Seasonality Index" (col P) values, drop the "Sales Volume Rank"
#   column (old Q) entirely and repurpose the old "Lifecycle Stage"
#   column (old R, now shifted to Q) with new per-week values.
# - Sheet "Summary": Max/Min Forecast Week become "N/A".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$weekStartDates = @{
    2  = "2025-02-02"
    3  = "2025-02-09"
    4  = "2025-02-16"
    5  = "2025-02-23"
    6  = "2025-03-02"
    7  = "2025-03-09"
    8  = "2025-03-16"
    9  = "2025-03-23"
    10 = "2025-03-30"
    11 = "2025-04-06"
    12 = "2025-04-13"
    13 = "2025-04-20"
    14 = "2025-04-27"
    15 = "2025-05-04"
    16 = "2025-05-11"
    17 = "2025-05-18"
}

$seasonalityIndex = @{
    2  = 1.1
    3  = 0.98
    4  = 1.17
    5  = 0.82
    6  = 1.16
    7  = 1.06
    8  = 1.02
    9  = 0.9
    10 = 0.82
    11 = 0.99
    12 = 1.06
    13 = 0.9
    14 = 0.86
    15 = 0.88
    16 = 1.07
    17 = 0.8
}

foreach ($row in 2..17) {
    # Week_Start_Date (column B) - force text so Excel doesn't reinterpret
    # the "YYYY-MM-DD" string as a date serial number.
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $weekStartDates[$row]
    $cellB.ClearFormats()

    # Sales Trend (column O)
    $ws.Cells.Item($row, 15).Value = "Stable"

    # Seasonality Index (column P)
    $ws.Cells.Item($row, 16).Value = $seasonalityIndex[$row]
}

# Drop the "Sales Volume Rank" column (Q). This shifts the old
# "Lifecycle Stage" column (R) left into Q, and shrinks the used range
# from A1:R17 to A1:Q17.
$ws.Range("Q1").EntireColumn.Delete()

# Update the (now shifted) Lifecycle Stage column (Q) with the new stage.
foreach ($row in 2..17) {
    $ws.Cells.Item($row, 17).Value = "Mature"
}

# Summary sheet: Max/Min Forecast Week -> "N/A"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(13, 2).Value = "N/A"
$ws2.Cells.Item(15, 2).Value = "N/A"
